$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C56").Formula = "'1"
Write-Host "Value: $($ws.Range('C56').Value())"
Write-Host "Type: $($ws.Range('C56').Value().GetType())"
Write-Host "NumberFormat: $($ws.Range('C56').NumberFormat())"
